$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete entire row 337 (the "眠たい目つきと暖かい微笑みで" post), shifting
# all subsequent rows up by one.
$ws.Rows.Item(337).Delete()
